$d = $word.ActiveDocument

function Get-ParaIndexByText($doc, $text) {
    $found = -1
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $pp = $doc.Paragraphs.Item($i)
        $t = $pp.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            $found = $i
        }
    }
    return $found
}

function Get-RangeStartOfText($doc, $text) {
    $r = $doc.Content
    $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $r.Start
}

# ---------------------------------------------------------------------------
# 1) "Window" + " Ratio Bounds" (two runs) -> a single run "Window Ratio Bounds"
#    Scope the Find/Replace to the Heading1 paragraph that holds this text so
#    the other occurrences of "Window Ratio" elsewhere in the body aren't
#    touched.
# ---------------------------------------------------------------------------
$headingIdx = Get-ParaIndexByText $d "Window Ratio Bounds"
$headingPara = $d.Paragraphs.Item($headingIdx)
$headingPara.Range.Find.Execute("Window Ratio Bounds", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "Window Ratio Bounds", 2)

# ---------------------------------------------------------------------------
# 2) Append a new "Layers" section (heading + three body paragraphs) right
#    after the "...fit the ratio." paragraph, before the trailing empty
#    paragraph / section break. The _GoBack bookmark that currently sits at
#    the end of that paragraph must end up at the end of the new last
#    paragraph ("To delete a layer...").
# ---------------------------------------------------------------------------

$openQuote = [char]0x201C
$closeQuote = [char]0x201D

$bodyPara2 = "After clicking the layers option at the top of the screen (image -> layers) a new window will appear where you can alter what layers there are in the program."
$bodyPara3 = "Adding a new layer is as simple as clicking " + $openQuote + "Add Layer" + $closeQuote + ". After that you can rename it by clicking " + $openQuote + "Rename" + $closeQuote + ". Layers can be moved up and down within the list of layers by clicking the " + $openQuote + "Move Up" + $closeQuote + " or " + $openQuote + "Move Down" + $closeQuote + " buttons. To select a layer to draw to, click on the layer name and then click " + $openQuote + "Select Layer" + $closeQuote + "."
$bodyPara4 = "To delete a layer, click on the layer to be deleted, and then click " + $openQuote + "Delete Layer" + $closeQuote

# Insert the final paragraph's text right before the _GoBack bookmark so the
# bookmark keeps tracking the end of the (new) last paragraph.
$bm = $d.Bookmarks("_GoBack")
$bm.Range.InsertBefore($bodyPara4)

# Split that text into its own paragraph, separate from "...fit the ratio."
$splitStart = Get-RangeStartOfText $d $bodyPara4
$d.Range($splitStart, $splitStart).InsertParagraphBefore()

# Insert the heading + the two other body paragraphs before that paragraph.
$insertStart = Get-RangeStartOfText $d $bodyPara4
$insertPoint = $d.Range($insertStart, $insertStart)
$insertPoint.InsertBefore("Layers`r$bodyPara2`r$bodyPara3`r")

# Apply the Heading1 style to the new "Layers" paragraph.
$layersIdx = Get-ParaIndexByText $d "Layers"
$layersPara = $d.Paragraphs.Item($layersIdx)
$layersPara.Style = "Heading 1"
